{"js": "// Edit the \"Nhi\u1ec7m v\u1ee5\" (Task) and \"K\u1ebft qu\u1ea3\" (Result) cells for the\n// \"\u0110o\u00e0n T\u00e2y \u0110\u00f4\" row of the status-report table:\n//   1. \"Hi\u1ec7n th\u1ef1c nh\u1eadn d\u1eef li\u1ec7u t\u1eeb Thingsboard.\"\n//        -> \"Hi\u1ec7n th\u1ef1c nh\u1eadn v\u00e0 hi\u1ec3n th\u1ecb d\u1eef li\u1ec7u t\u1eeb Thingsboard.\"\n//   2. \"\u0110ang hi\u1ec7n th\u1ef1c h\u00e0m nh\u1eadn d\u1eef li\u1ec7u t\u1eeb MQTT. G\u1eb7p v\u1ea5n \u0111\u1ec1 do kh\u00f4ng th\u1ec3\n//      2 thi\u1ebft b\u1ecb k\u1ebft n\u1ed1i Thingsboard c\u00f9ng l\u00fac.\"\n//        -> \"\u0110\u00e3 t\u1ea1o ph\u1ea7n hi\u1ec3n th\u1ecb d\u1eef li\u1ec7u tr\u00ean Unity. Tuy nhi\u00ean, code \u0111\u1ec3\n//      nh\u1eadn d\u1eef li\u1ec7u \u0111\u00e3 k\u1ebft n\u1ed1i \u0111\u01b0\u1ee3c v\u1edbi feed tr\u00ean adafruit nh\u01b0ng l\u1ea1i\n//      kh\u00f4ng g\u1eedi d\u1eef li\u1ec7u v\u1ec1.\"\n\nconst body = context.document.body;\n\n// --- Edit 1: insert \" v\u00e0 hi\u1ec3n th\u1ecb\" right after \"Hi\u1ec7n th\u1ef1c nh\u1eadn\" -----------\nconst firstSearch = body.search(\"Hi\u1ec7n th\u1ef1c nh\u1eadn d\u1eef li\u1ec7u t\u1eeb Thingsboard.\", {\n  matchCase: true,\n});\nfirstSearch.load(\"items\");\nawait context.sync();\n\nif (firstSearch.items.length > 0) {\n  const target = firstSearch.items[0];\n  target.insertText(\n    \"Hi\u1ec7n th\u1ef1c nh\u1eadn v\u00e0 hi\u1ec3n th\u1ecb d\u1eef li\u1ec7u t\u1eeb Thingsboard.\",\n    \"Replace\"\n  );\n  await context.sync();\n} else {\n  // Fallback: the text may have already been partially edited / split into\n  // runs, so locate just the leading fragment and splice the addition in.\n  const partial = body.search(\"Hi\u1ec7n th\u1ef1c nh\u1eadn\", { matchCase: true });\n  partial.load(\"items\");\n  await context.sync();\n  if (partial.items.length > 0) {\n    partial.items[0].insertText(\" v\u00e0 hi\u1ec3n th\u1ecb\", \"After\");\n    await context.sync();\n  }\n}\n\n// --- Edit 2: replace the whole \"K\u1ebft qu\u1ea3\" text for that row -----------------\nconst secondSearch = body.search(\n  \"\u0110ang hi\u1ec7n th\u1ef1c h\u00e0m nh\u1eadn d\u1eef li\u1ec7u t\u1eeb MQTT. G\u1eb7p v\u1ea5n \u0111\u1ec1 do kh\u00f4ng th\u1ec3 2 thi\u1ebft b\u1ecb k\u1ebft n\u1ed1i Thingsboard c\u00f9ng l\u00fac.\",\n  { matchCase: true }\n);\nsecondSearch.load(\"items\");\nawait context.sync();\n\nif (secondSearch.items.length > 0) {\n  secondSearch.items[0].insertText(\n    \"\u0110\u00e3 t\u1ea1o ph\u1ea7n hi\u1ec3n th\u1ecb d\u1eef li\u1ec7u tr\u00ean Unity. Tuy nhi\u00ean, code \u0111\u1ec3 nh\u1eadn d\u1eef li\u1ec7u \u0111\u00e3 k\u1ebft n\u1ed1i \u0111\u01b0\u1ee3c v\u1edbi feed tr\u00ean adafruit nh\u01b0ng l\u1ea1i kh\u00f4ng g\u1eedi d\u1eef li\u1ec7u v\u1ec1.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "# Edit the \"Nhi\u1ec7m v\u1ee5\" (Task) and \"K\u1ebft qu\u1ea3\" (Result) cells for the\n# \"\u0110o\u00e0n T\u00e2y \u0110\u00f4\" row of the status-report table:\n#   1. \"Hi\u1ec7n th\u1ef1c nh\u1eadn d\u1eef li\u1ec7u t\u1eeb Thingsboard.\"\n#        -> \"Hi\u1ec7n th\u1ef1c nh\u1eadn v\u00e0 hi\u1ec3n th\u1ecb d\u1eef li\u1ec7u t\u1eeb Thingsboard.\"\n#   2. \"\u0110ang hi\u1ec7n th\u1ef1c h\u00e0m nh\u1eadn d\u1eef li\u1ec7u t\u1eeb MQTT. G\u1eb7p v\u1ea5n \u0111\u1ec1 do kh\u00f4ng th\u1ec3\n#      2 thi\u1ebft b\u1ecb k\u1ebft n\u1ed1i Thingsboard c\u00f9ng l\u00fac.\"\n#        -> \"\u0110\u00e3 t\u1ea1o ph\u1ea7n hi\u1ec3n th\u1ecb d\u1eef li\u1ec7u tr\u00ean Unity. Tuy nhi\u00ean, code \u0111\u1ec3\n#      nh\u1eadn d\u1eef li\u1ec7u \u0111\u00e3 k\u1ebft n\u1ed1i \u0111\u01b0\u1ee3c v\u1edbi feed tr\u00ean adafruit nh\u01b0ng l\u1ea1i\n#      kh\u00f4ng g\u1eedi d\u1eef li\u1ec7u v\u1ec1.\"\n\n$d = $word.ActiveDocument\n\n$oldTask = \"Hi\u1ec7n th\u1ef1c nh\u1eadn d\u1eef li\u1ec7u t\u1eeb Thingsboard.\"\n$newTask = \"Hi\u1ec7n th\u1ef1c nh\u1eadn v\u00e0 hi\u1ec3n th\u1ecb d\u1eef li\u1ec7u t\u1eeb Thingsboard.\"\n\n$oldResult = \"\u0110ang hi\u1ec7n th\u1ef1c h\u00e0m nh\u1eadn d\u1eef li\u1ec7u t\u1eeb MQTT. G\u1eb7p v\u1ea5n \u0111\u1ec1 do kh\u00f4ng th\u1ec3 2 thi\u1ebft b\u1ecb k\u1ebft n\u1ed1i Thingsboard c\u00f9ng l\u00fac.\"\n$newResult = \"\u0110\u00e3 t\u1ea1o ph\u1ea7n hi\u1ec3n th\u1ecb d\u1eef li\u1ec7u tr\u00ean Unity. Tuy nhi\u00ean, code \u0111\u1ec3 nh\u1eadn d\u1eef li\u1ec7u \u0111\u00e3 k\u1ebft n\u1ed1i \u0111\u01b0\u1ee3c v\u1edbi feed tr\u00ean adafruit nh\u01b0ng l\u1ea1i kh\u00f4ng g\u1eedi d\u1eef li\u1ec7u v\u1ec1.\"\n\n$targetCell = $null\n$resultCell = $null\n\n$table = $d.Tables(1)\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n    $taskCell = $table.Cell($r, 2)\n    $taskText = $taskCell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($taskText -eq $oldTask) {\n        $targetCell = $taskCell\n        $resultCell = $table.Cell($r, 3)\n        break\n    }\n}\n\nif ($targetCell -ne $null) {\n    $rng = $targetCell.Range\n    $rng.End = $rng.End - 1\n    $rng.Text = $newTask\n\n    $rng2 = $resultCell.Range\n    $rng2.End = $rng2.End - 1\n    $rng2.Text = $newResult\n} else {\n    # Fallback: use Find/Replace in case the table lookup above did not\n    # locate the row (e.g. text already partially edited).\n    $find1 = $d.Content.Find\n    $find1.ClearFormatting()\n    $find1.Text = $oldTask\n    $find1.Replacement.ClearFormatting()\n    $find1.Replacement.Text = $newTask\n    $find1.Execute([ref]$oldTask, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newTask, 2) | Out-Null\n\n    $find2 = $d.Content.Find\n    $find2.ClearFormatting()\n    $find2.Text = $oldResult\n    $find2.Replacement.ClearFormatting()\n    $find2.Replacement.Text = $newResult\n    $find2.Execute([ref]$oldResult, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newResult, 2) | Out-Null\n}\n"}
